$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume-change (E) figures.
# D-column values are stored as text in the workbook (e.g. "7.00", "1.00",
# "47.801.18"); several of the new values would otherwise be auto-parsed by
# Excel as numbers (losing the trailing zeros / thousands dots), so we force
# those cells to Text format before writing them.

$ws.Range("D2").Value = '47.801.18'
$ws.Range("E2").Value = '  -0.88%  '
$ws.Range("D3").Value = '2.477.37'
$ws.Range("E3").Value = '  -1.81%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.88'
$ws.Range("E5").Value = '  -2.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.72'
$ws.Range("E6").Value = '  -5.25%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.515'
$ws.Range("E7").Value = '  -3.25%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.531'
$ws.Range("E9").Value = '  -3.96%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.56'
$ws.Range("E10").Value = '  -4.94%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.34'
$ws.Range("E11").Value = '  -1.15%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0795'
$ws.Range("E12").Value = '  -3.40%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.00'
$ws.Range("E14").Value = '  -3.90%  '
$ws.Range("D15").Value = '2.866.40'
$ws.Range("D16").Value = '2.492.61'
$ws.Range("E16").Value = '  -1.40%  '
$ws.Range("E17").Value = '  -4.09%  '
$ws.Range("D18").Value = '47.726.10'
$ws.Range("E18").Value = '  -0.68%  '
$ws.Range("E19").Value = '  +7.64%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.57'
$ws.Range("E20").Value = '  -6.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.49'
$ws.Range("E21").Value = '  -2.20%  '
$ws.Range("D22").Value = '0.0₃0923'
$ws.Range("E22").Value = '  -2.61%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '280.10'
$ws.Range("E23").Value = '  +5.81%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.70'
$ws.Range("E24").Value = '  -1.83%  '
$ws.Range("E25").Value = '  -3.51%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.998'
$ws.Range("E26").Value = '  +0.12%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.66'
$ws.Range("E27").Value = '  -1.60%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.22'
$ws.Range("E28").Value = '  +0.52%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.52'
$ws.Range("E29").Value = '  -5.89%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.135'
$ws.Range("E30").Value = '  -5.38%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.31'
$ws.Range("E31").Value = '  -4.76%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.20'
$ws.Range("E32").Value = '  -1.21%  '
$ws.Range("E33").Value = '  -0.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.79'
$ws.Range("E34").Value = '  -5.31%  '
$ws.Range("E35").Value = '  -3.57%  '
$ws.Range("E36").Value = '  -3.30%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.92'
$ws.Range("E37").Value = '  -3.65%  '
$ws.Range("E38").Value = '  -5.78%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.84'
$ws.Range("E39").Value = '  -5.16%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '122.45'
$ws.Range("E40").Value = '  +1.30%  '
$ws.Range("E41").Value = '  -2.09%  '
$ws.Range("E42").Value = '  -0.17%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.82'
$ws.Range("E43").Value = '  -1.37%  '
$ws.Range("E44").Value = '  -1.27%  '
$ws.Range("D45").Value = '1.987.27'
$ws.Range("E45").Value = '  -1.58%  '
$ws.Range("E46").Value = '  -1.88%  '
$ws.Range("E47").Value = '  -2.17%  '
$ws.Range("E48").Value = '  -3.45%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.87'
$ws.Range("E49").Value = '  -2.98%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.07'
$ws.Range("E50").Value = '  -3.03%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.86'
$ws.Range("E51").Value = '  -0.41%  '
